$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.325.56"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.667.20"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5304"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2646"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06364"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.531"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "1.674.52"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "1.896.36"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5604"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").Value = "0.0₅8140"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "26.339.83"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "198.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.84%  "
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.060"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1213"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.243"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05885"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.285"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.547"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.323"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.832"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9617"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.432"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5811"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.965"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("D41").Value = "1.074.05"
$ws.Range("E41").Value = "  +3.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8562"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").Value = "1.807.63"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.015"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈105"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4412"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.054"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05146"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
